$d = $word.ActiveDocument

# The list item currently reads (runs joined):
#   "https://2phimmoi.net/ - Tải được (chất lượng kém, Huong.Dan.Lay.Phim.2Phimmoi.docx)."
# and must become:
#   "https://2phimmoi.net/ - Tải được bằng IDM (chất lượng kém)."
$marker = "Huong.Dan.Lay.Phim.2Phimmoi.docx"

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$marker*") {
        $pRange = $p.Range

        # Rebuild the paragraph with the new run split described by the
        # change: " được (" + "chất lượng kém, " + "Huong...docx" (3 runs)
        # becomes " được" + " bằng IDM" + " (" + "chất lượng kém" (4 runs),
        # while the leading runs (hyperlink, " - ", "Tải") and the trailing
        # run (").") are kept byte-for-byte as they were.
        $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7E0BE5A6" w14:textId="4A458464" w:rsidR="000E352D" w:rsidRDefault="00BB28C3" w:rsidP="00BB28C3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:hyperlink r:id="rId5" w:history="1"><w:r w:rsidRPr="00AE3527"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://2phimmoi.net/</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve"> - </w:t></w:r><w:r w:rsidR="00267EAF"><w:t>Tải</w:t></w:r><w:r><w:t xml:space="preserve"> được</w:t></w:r><w:r><w:t xml:space="preserve"> bằng IDM</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:t>chất lượng kém</w:t></w:r><w:r><w:t>).</w:t></w:r></w:p>'

        $pRange.InsertXML($newParaXml)
        break
    }
}

# InsertXML re-serialises the whole paragraph, which drops the rStyle
# reference inside the hyperlink run's <w:rPr> (a quirk of this host).
# Restore the "Hyperlink" character style on that run so the link keeps
# its original appearance.
$linkRange = $d.Content
if ($linkRange.Find.Execute("https://2phimmoi.net/", $false, $false, $false, $false, $false, $true, 1, $false)) {
    $linkRange.Style = "Hyperlink"
}
